$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hot fix for MS1 centroid Thermo .mzML files in Shotgun mode:
# swap the H2O flags for the sn1/sn2 ([M-H]-sn1, [M-H]-sn2) rows
# with the sn1-H2O/sn2-H2O rows.
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0

# Move the active selection to F7 to match the saved cursor position.
$ws.Range("F7").Select()
